# Scheduled-runner price/profit refresh for the Famfrit_Profits sheets.
# Updates the market-price-derived columns (H:N) on a handful of rows
# across all eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 7018.875
$ws.Cells.Item(9, 9).Value = 10371.5
$ws.Cells.Item(9, 11).Value = 10371.5
$ws.Cells.Item(9, 13).Value = -10202.5
$ws.Cells.Item(32, 8).Value = 7152.4287
$ws.Cells.Item(32, 10).Value = 7515.6
$ws.Cells.Item(32, 12).Value = 7515.6
$ws.Cells.Item(32, 14).Value = -8167.6
$ws.Cells.Item(62, 8).Value = 5824.2144
$ws.Cells.Item(62, 9).Value = 4309
$ws.Cells.Item(62, 10).Value = 6666
$ws.Cells.Item(62, 11).Value = 4309
$ws.Cells.Item(62, 12).Value = 6666
$ws.Cells.Item(62, 13).Value = -3685
$ws.Cells.Item(62, 14).Value = -7914
$ws.Cells.Item(65, 8).Value = 5824.2144
$ws.Cells.Item(65, 9).Value = 4309
$ws.Cells.Item(65, 10).Value = 6666
$ws.Cells.Item(65, 11).Value = 21545
$ws.Cells.Item(65, 12).Value = 33330
$ws.Cells.Item(65, 13).Value = -18425
$ws.Cells.Item(65, 14).Value = -39570
$ws.Cells.Item(99, 8).Value = 695.1818
$ws.Cells.Item(99, 9).Value = 664.7
$ws.Cells.Item(99, 10).Value = 1000
$ws.Cells.Item(99, 11).Value = 1994.1
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = -496.1000000000001
$ws.Cells.Item(99, 14).Value = -5996
$ws.Cells.Item(135, 8).Value = 14707067
$ws.Cells.Item(135, 9).Value = 867
$ws.Cells.Item(135, 11).Value = 7803
$ws.Cells.Item(135, 13).Value = -5268
$ws.Cells.Item(137, 8).Value = 2339.9768
$ws.Cells.Item(137, 9).Value = 2398.1155
$ws.Cells.Item(137, 10).Value = 2251.0588
$ws.Cells.Item(137, 11).Value = 7194.3465
$ws.Cells.Item(137, 12).Value = 6753.176399999999
$ws.Cells.Item(137, 13).Value = -4644.3465
$ws.Cells.Item(137, 14).Value = -11853.1764

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15390669
$ws.Cells.Item(32, 9).Value = 18521836
$ws.Cells.Item(32, 10).Value = 19481.637
$ws.Cells.Item(32, 11).Value = 18521836
$ws.Cells.Item(32, 12).Value = 19481.637
$ws.Cells.Item(32, 13).Value = -18521549
$ws.Cells.Item(32, 14).Value = -20055.637
$ws.Cells.Item(45, 8).Value = 2117.6
$ws.Cells.Item(45, 9).Value = 1893.5
$ws.Cells.Item(45, 11).Value = 1893.5
$ws.Cells.Item(45, 13).Value = -1516.5
$ws.Cells.Item(97, 8).Value = 1549.9615
$ws.Cells.Item(97, 9).Value = 1622.9048
$ws.Cells.Item(97, 11).Value = 1622.9048
$ws.Cells.Item(97, 13).Value = -1126.9048

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 4064.7058
$ws.Cells.Item(107, 9).Value = 3617.25
$ws.Cells.Item(107, 11).Value = 3617.25
$ws.Cells.Item(107, 13).Value = -1697.25
$ws.Cells.Item(132, 8).Value = 120000
$ws.Cells.Item(132, 10).Value = 120000
$ws.Cells.Item(132, 12).Value = 120000
$ws.Cells.Item(132, 14).Value = -130120

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 1876.75
$ws.Cells.Item(122, 9).Value = 1807.3636
$ws.Cells.Item(122, 11).Value = 5422.0908
$ws.Cells.Item(122, 13).Value = -2972.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 942.6
$ws.Cells.Item(121, 9).Value = 173.625
$ws.Cells.Item(121, 10).Value = 1821.4286
$ws.Cells.Item(121, 11).Value = 520.875
$ws.Cells.Item(121, 12).Value = 5464.2858
$ws.Cells.Item(121, 13).Value = 789.125
$ws.Cells.Item(121, 14).Value = -8084.2858
$ws.Cells.Item(122, 8).Value = 900
$ws.Cells.Item(122, 9).Value = 975
$ws.Cells.Item(122, 10).Value = 750
$ws.Cells.Item(122, 11).Value = 8775
$ws.Cells.Item(122, 12).Value = 6750
$ws.Cells.Item(122, 13).Value = -6325
$ws.Cells.Item(122, 14).Value = -11650
$ws.Cells.Item(131, 8).Value = 1603.8
$ws.Cells.Item(131, 10).Value = 1824.3334
$ws.Cells.Item(131, 12).Value = 5473.0002
$ws.Cells.Item(131, 14).Value = -15553.0002
$ws.Cells.Item(133, 8).Value = 10114.25
$ws.Cells.Item(133, 10).Value = 15126.5
$ws.Cells.Item(133, 12).Value = 45379.5
$ws.Cells.Item(133, 14).Value = -55499.5
$ws.Cells.Item(134, 8).Value = 3600.2632
$ws.Cells.Item(134, 9).Value = 1583.8125
$ws.Cells.Item(134, 11).Value = 4751.4375
$ws.Cells.Item(134, 13).Value = 318.5625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3763.0625
$ws.Cells.Item(113, 9).Value = 2775.875
$ws.Cells.Item(113, 11).Value = 2775.875
$ws.Cells.Item(113, 13).Value = -605.875
$ws.Cells.Item(117, 8).Value = 49000
$ws.Cells.Item(117, 10).Value = 49000
$ws.Cells.Item(117, 12).Value = 49000
$ws.Cells.Item(117, 14).Value = -55884
$ws.Cells.Item(118, 8).Value = 20060.334
$ws.Cells.Item(118, 10).Value = 19900
$ws.Cells.Item(118, 12).Value = 19900
$ws.Cells.Item(118, 14).Value = -23214

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4236.2607
$ws.Cells.Item(7, 9).Value = 4388
$ws.Cells.Item(7, 10).Value = 4000.2222
$ws.Cells.Item(7, 11).Value = 4388
$ws.Cells.Item(7, 12).Value = 4000.2222
$ws.Cells.Item(7, 13).Value = -4276
$ws.Cells.Item(7, 14).Value = -4224.2222
$ws.Cells.Item(46, 8).Value = 1344.8723
$ws.Cells.Item(46, 10).Value = 3078.6155
$ws.Cells.Item(46, 12).Value = 3078.6155
$ws.Cells.Item(46, 14).Value = -3454.6155
$ws.Cells.Item(93, 8).Value = 2635.9443
$ws.Cells.Item(93, 9).Value = 1168.5
$ws.Cells.Item(93, 10).Value = 3809.9
$ws.Cells.Item(93, 11).Value = 1168.5
$ws.Cells.Item(93, 12).Value = 3809.9
$ws.Cells.Item(93, 13).Value = 79.5
$ws.Cells.Item(93, 14).Value = -6305.9
$ws.Cells.Item(100, 8).Value = 5032.8335
$ws.Cells.Item(100, 9).Value = 3499.5
$ws.Cells.Item(100, 10).Value = 5799.5
$ws.Cells.Item(100, 11).Value = 3499.5
$ws.Cells.Item(100, 12).Value = 5799.5
$ws.Cells.Item(100, 13).Value = -2958.5
$ws.Cells.Item(100, 14).Value = -6881.5
$ws.Cells.Item(126, 8).Value = 4236.2607
$ws.Cells.Item(126, 9).Value = 4388
$ws.Cells.Item(126, 10).Value = 4000.2222
$ws.Cells.Item(126, 11).Value = 13164
$ws.Cells.Item(126, 12).Value = 12000.6666
$ws.Cells.Item(126, 13).Value = -10694
$ws.Cells.Item(126, 14).Value = -16940.6666
$ws.Cells.Item(132, 8).Value = 95240950
$ws.Cells.Item(132, 9).Value = 2809.0715
$ws.Cells.Item(132, 11).Value = 8427.2145
$ws.Cells.Item(132, 13).Value = -5897.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(11, 8).Value = 3000
$ws.Cells.Item(11, 10).Value = 3000
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 14).Value = -3284
$ws.Cells.Item(81, 8).Value = 499.77777
$ws.Cells.Item(81, 9).Value = 437.25
$ws.Cells.Item(81, 10).Value = 1000
$ws.Cells.Item(81, 11).Value = 874.5
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 13).Value = 186.5
$ws.Cells.Item(81, 14).Value = -4122
$ws.Cells.Item(84, 8).Value = 499.77777
$ws.Cells.Item(84, 9).Value = 437.25
$ws.Cells.Item(84, 10).Value = 1000
$ws.Cells.Item(84, 11).Value = 4372.5
$ws.Cells.Item(84, 12).Value = 10000
$ws.Cells.Item(84, 13).Value = 931.5
$ws.Cells.Item(84, 14).Value = -4122
